{"js": "// Replace each \"old\" three-digit x one-digit multiplication fact with its \"new\" counterpart.\n// Pairs taken 1:1, in document order, from the authoritative diff.\nconst replacements = [\n  [\"371\u00d72=742\", \"517\u00d73=1551\"],\n  [\"743\u00d72=1486\", \"626\u00d78=5008\"],\n  [\"796\u00d74=3184\", \"648\u00d78=5184\"],\n  [\"131\u00d75=655\", \"888\u00d79=7992\"],\n  [\"693\u00d72=1386\", \"197\u00d77=1379\"],\n  [\"398\u00d73=1194\", \"333\u00d76=1998\"],\n  [\"429\u00d74=1716\", \"803\u00d76=4818\"],\n  [\"925\u00d75=4625\", \"722\u00d75=3610\"],\n  [\"622\u00d77=4354\", \"779\u00d73=2337\"],\n  [\"826\u00d79=7434\", \"498\u00d75=2490\"],\n  [\"838\u00d76=5028\", \"161\u00d75=805\"],\n  [\"455\u00d77=3185\", \"680\u00d78=5440\"],\n  [\"171\u00d74=684\", \"612\u00d78=4896\"],\n  [\"809\u00d72=1618\", \"237\u00d75=1185\"],\n  [\"514\u00d76=3084\", \"356\u00d76=2136\"],\n  [\"210\u00d72=420\", \"252\u00d77=1764\"],\n  [\"712\u00d73=2136\", \"164\u00d74=656\"],\n  [\"358\u00d79=3222\", \"503\u00d74=2012\"],\n  [\"494\u00d73=1482\", \"353\u00d79=3177\"],\n  [\"474\u00d75=2370\", \"884\u00d73=2652\"],\n  [\"925\u00d74=3700\", \"233\u00d76=1398\"],\n  [\"774\u00d78=6192\", \"728\u00d73=2184\"],\n  [\"140\u00d79=1260\", \"125\u00d73=375\"],\n  [\"595\u00d74=2380\", \"782\u00d74=3128\"],\n  [\"984\u00d76=5904\", \"631\u00d73=1893\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  // matchCase keeps the search from colliding across distinct-but-overlapping digit runs.\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"old\" three-digit x one-digit multiplication fact with its \"new\" counterpart.\n# Pairs taken 1:1, in document order, from the authoritative diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"371\u00d72=742\", \"517\u00d73=1551\"),\n  @(\"743\u00d72=1486\", \"626\u00d78=5008\"),\n  @(\"796\u00d74=3184\", \"648\u00d78=5184\"),\n  @(\"131\u00d75=655\", \"888\u00d79=7992\"),\n  @(\"693\u00d72=1386\", \"197\u00d77=1379\"),\n  @(\"398\u00d73=1194\", \"333\u00d76=1998\"),\n  @(\"429\u00d74=1716\", \"803\u00d76=4818\"),\n  @(\"925\u00d75=4625\", \"722\u00d75=3610\"),\n  @(\"622\u00d77=4354\", \"779\u00d73=2337\"),\n  @(\"826\u00d79=7434\", \"498\u00d75=2490\"),\n  @(\"838\u00d76=5028\", \"161\u00d75=805\"),\n  @(\"455\u00d77=3185\", \"680\u00d78=5440\"),\n  @(\"171\u00d74=684\", \"612\u00d78=4896\"),\n  @(\"809\u00d72=1618\", \"237\u00d75=1185\"),\n  @(\"514\u00d76=3084\", \"356\u00d76=2136\"),\n  @(\"210\u00d72=420\", \"252\u00d77=1764\"),\n  @(\"712\u00d73=2136\", \"164\u00d74=656\"),\n  @(\"358\u00d79=3222\", \"503\u00d74=2012\"),\n  @(\"494\u00d73=1482\", \"353\u00d79=3177\"),\n  @(\"474\u00d75=2370\", \"884\u00d73=2652\"),\n  @(\"925\u00d74=3700\", \"233\u00d76=1398\"),\n  @(\"774\u00d78=6192\", \"728\u00d73=2184\"),\n  @(\"140\u00d79=1260\", \"125\u00d73=375\"),\n  @(\"595\u00d74=2380\", \"782\u00d74=3128\"),\n  @(\"984\u00d76=5904\", \"631\u00d73=1893\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  # 2 = wdReplaceAll (replace every exact-case match of $oldText with $newText).\n  $found = $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n  if (-not $found) {\n    throw \"No match found for: $oldText\"\n  }\n}\n"}
